$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 21 de Marzo de 2020 a las 14:16"

# Provincia/city names, in row order (rows 4-61), after re-sorting by Casos totales desc
$cities = @(
    "Madrid",
    "Cataluña",
    "Araba/Alava",
    "Navarra",
    "Valencia/Valencia",
    "Bizkaia/Vizcaya",
    "La Rioja",
    "Asturias",
    "Malaga",
    "Ciudad Real",
    "Alacant/Alicante",
    "Toledo",
    "Albacete",
    "Zaragoza",
    "A Coruña",
    "Burgos",
    "Salamanca",
    "Illes Balears",
    "Granada",
    "Caceres",
    "Valladolid",
    "Guadalajara",
    "Gipuzkoa/Guipuzcoa",
    "Murcia",
    "Cantabria",
    "Leon",
    "Pontevedra",
    "Tenerife",
    "Sevilla",
    "Aragon",
    "Illes Balears*",
    "Segovia",
    "Cordoba",
    "Jaen",
    "Soria",
    "Avila",
    "Badajoz",
    "Castello/Castellon",
    "Cadiz",
    "Cuenca",
    "Gran Canaria",
    "Ourense",
    "Zamora",
    "Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena",
    "Almeria",
    "Lugo",
    "Palencia",
    "Teruel",
    "Huelva",
    "Huesca",
    "Melilla",
    "Fuerteventura",
    "Arroyo de la Luz",
    "La Palma",
    "Ceuta",
    "Lanzarote",
    "La Gomera",
    "El Hierro"
)

# Casos totales, Casos activos, Recuperados, Muertes (parallel arrays, rows 4-61)
$casosTotales = @(8921, 4203, 703, 664, 627, 539, 497, 486, 424, 400, 372, 370, 327, 278, 270, 269, 265, 246, 244, 243, 241, 237, 223, 215, 215, 201, 193, 192, 190, 174, 169, 157, 125, 119, 119, 114, 111, 104, 103, 94, 70, 63, 59, 58, 49, 46, 41, 40, 36, 34, 25, 12, 7, 7, 5, 3, 3, 1)
$casosActivos = @(1186, 3, 21, 2, 12, 21, 13, 12, 72, 8, 7, 15, 8, 0, 5, 27, 13, 10, 72, 2, 13, 2, 21, 1, 11, 3, 5, 4, 72, 0, 6, 9, 72, 72, 5, 14, 5, 1, 72, 5, 0, 5, 3, 0, 72, 5, 2, 0, 72, 0, 0, 0, 0, 0, 0, 0, 2, 0)
$recuperados  = @(5351, 4078, 655, 582, 600, 522, 564, 469, 406, 364, 348, 336, 291, 264, 267, 175, 180, 189, 234, 194, 193, 231, 217, 213, 132, 156, 191, 184, 188, 163, 161, 140, 123, 117, 71, 55, 86, 102, 102, 84, 69, 63, 42, 58, 49, 45, 25, 38, 36, 34, 25, 12, 7, 7, 5, 3, 1, 1)
$muertes      = @(804, 122, 48, 10, 15, 17, 15, 7, 18, 28, 17, 19, 28, 14, 3, 14, 15, 4, 10, 10, 10, 4, 6, 1, 4, 7, 2, 4, 1, 11, 2, 13, 2, 2, 6, 7, 2, 1, 1, 5, 1, 0, 2, 3, 0, 1, 0, 2, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)

$startRow = 4
for ($i = 0; $i -lt $cities.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $cities[$i]
    $ws.Cells.Item($r, 2).Value = $casosTotales[$i]
    $ws.Cells.Item($r, 3).Value = $casosActivos[$i]
    $ws.Cells.Item($r, 4).Value = $recuperados[$i]
    $ws.Cells.Item($r, 5).Value = $muertes[$i]
}
